# Applies the "Add files via upload" revision to the Cohere evaluation workbook.
# The edit corrects several machine-translated Telugu strings (column C,
# "Cohere_Translation") together with their associated BLEU scores
# (column D), and re-saves the workbook the way the file was uploaded
# (worksheet renamed to the default "Sheet1", cursor left on C143).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet back to the generic default name ---
$ws.Name = "Sheet1"

# --- Corrected translations + recomputed BLEU scores ---

# Row 6
$ws.Range("C6").Value = "సింహాల గుంపులు పాటిస్తారు కుక్కల కూడలో అథవా కొక్కల కూడలో ఉండటం మాత్రమే కాదు, వాటి వ్యవహారం సింహాలతో సంతోషంగా సమానం (కింద పెద్ద పందులతో కాదు) మరియు వాటి బాధిత ప్రాణికి అత్యంత మరణాన్నికంగా ఉండటం."
$ws.Range("D6").Value = 0.01435802452844482

# Row 77
$ws.Range("C77").Value = "ప్రత్యేకించి, అందుబాటులో వ్యక్తి దుష్టువ్వాడా లేదా సత్యం వినియోగిస్తానను సరిగ్గా అర్థం చేసి అందుబాటులను కనుగొనబడుతుంది అని అవగాహన చేస్తారు."
$ws.Range("D77").Value = 0.01973065769362093

# Row 89
$ws.Range("C89").Value = "విదేశాన్ని సంబంధించిన లింగ పునర్నిర్ణయ శస్త్రచికిక పై యోజన చేయునప్పుళ్ళిని వారు తిరిగి వెళ్లడానికి సరిపోయిన దాఖలు తీసుకోవాలి ఉంటుంది."
$ws.Range("D89").Value = 0.01659143932516396

# Row 99 (BLEU score unchanged)
$ws.Range("C99").Value = "ఫాస్టర్ కేర్ అన్ని అవసరాలను పూర్తి చేయాలని ఉంది వాటిని వాటి ముందుగా ఉన్న వీటిలో కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్ని కొన్"

# Row 100
$ws.Range("C100").Value = "చెప్పబడిన పని సిద్ధాంతాత్మకంగా ఉంది, కాబట్టికి సిమ్యులేటర్ వలన సాగిటారియస్ గ్యాలక్సీ గురించి వివరాలను సిరిసింపాడిన ప్రోగ్రాంతను రచించారు."
$ws.Range("D100").Value = 0.01947104482106315

# Row 107
$ws.Range("C107").Value = "మోసాసారస్ ప్రాచీన కాలంలో అతిపెద్ద అంత్యాయినియామికి మాత్రమే భయం లేదు."
$ws.Range("D107").Value = 0.01350950256839151

# Row 131
$ws.Range("C131").Value = "ఈ పరిణామంగా, రెండు చేతి జాతీయ పశువులు అంత్యమవుతున్నాయి మరియు రెండు చేతి చేతులు అపాయకరమైనవిగా మారివుతున్నాయి, దానికి సంబంధించిన గుండెముక్క చేతి హంప్‌బ్యాక్ చూబ్ అని అనువదించబడింది."
$ws.Range("D131").Value = 0.01251837734451281

# Row 143
$ws.Range("C143").Value = "దువల్ దంపతిరాలారు మరియు రెండు పెద్ద పిల్లలు ఉన్నారు. కాని, కథను చెప్పిన మిల్లర్ మీద అతనికి చాలా ప్రభావం చేసినది లేదు."
$ws.Range("D143").Value = 0.0149361928534572

# Row 149
$ws.Range("C149").Value = "పోలీస్ వివరాల ప్రకారం, ఫోటోగ్రాఫర్ని తగ్గిన వాహనాన్ని డ్రైవ్ చేసిన వ్యక్తి క్రిమినల్ కేసుల ఎదురు అసమర్థను అనుభవించకూడదు."
$ws.Range("D149").Value = 0.04126321688145509

# --- Leave the view the way it was when the file was saved ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 139
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("C143").Select()
